$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.243.59"
$ws.Range("E2").Value = "  +0.80%  "

$ws.Range("D3").Value = "1.872.83"
$ws.Range("E3").Value = "  +3.91%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.68"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5015"
$ws.Range("E7").Value = "  -1.22%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3923"
$ws.Range("E8").Value = "  +1.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.09583"
$ws.Range("E9").Value = "  +7.60%  "

$ws.Range("E10").Value = "  +4.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "40.99"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.484"
$ws.Range("E12").Value = "  +1.62%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.01"
$ws.Range("E13").Value = "  +3.34%  "

$ws.Range("D14").Value = "1.864.45"
$ws.Range("E14").Value = "  +3.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.000"
$ws.Range("E15").Value = "  +0.03%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.414"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001127"
$ws.Range("E17").Value = "  +1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.19"
$ws.Range("E18").Value = "  +1.18%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06631"
$ws.Range("E19").Value = "  +0.76%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.49"
$ws.Range("E20").Value = "  +1.30%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.136"
$ws.Range("E22").Value = "  +2.02%  "

$ws.Range("D23").Value = "28.291.58"
$ws.Range("E23").Value = "  +0.84%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.34"
$ws.Range("E24").Value = "  +2.59%  "

$ws.Range("E25").Value = "  +2.89%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.544"
$ws.Range("E26").Value = "  +5.70%  "

$ws.Range("D27").Value = "2.082.73"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "21.19"
$ws.Range("E28").Value = "  +4.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "157.74"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.35"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.067"
$ws.Range("E31").Value = "  +2.03%  "

$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1055"
$ws.Range("E32").Value = "  -3.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.641"
$ws.Range("E33").Value = "  +1.41%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.627"
$ws.Range("E34").Value = "  -0.33%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.06756"
$ws.Range("E35").Value = "  -2.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.497"
$ws.Range("E36").Value = "  +5.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02397"
$ws.Range("E37").Value = "  +2.63%  "

$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.50"
$ws.Range("E39").Value = "  +0.78%  "

$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6355"
$ws.Range("E40").Value = "  +3.77%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.974"
$ws.Range("E41").Value = "  -0.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.179"
$ws.Range("E42").Value = "  +2.47%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.002"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.59"
$ws.Range("E44").Value = "  +2.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6055"
$ws.Range("E45").Value = "  +2.65%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.665"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.264"
$ws.Range("E47").Value = "  -2.29%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "123.89"
$ws.Range("E48").Value = "  -1.11%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.989"
$ws.Range("E49").Value = "  +2.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.195"
$ws.Range("E50").Value = "  +0.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06837"
$ws.Range("E51").Value = "  +1.42%  "

